# Weekly update of fruit/vegetable (hortaliza) price data.
# The historical rows 2-5 are re-ordered: the data that used to sit on
# row 2 now belongs on row 4, row 3's data moves to row 5, row 4's data
# moves to row 2, and row 5's data moves to row 3 (i.e. row2<->row5 and
# row3<->row4 swap their contents for columns D,I,J,K,L,M,N,P,Q).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that vary between rows (date, quality, volume, prices, unit, etc.)
$cols = @("D", "I", "J", "K", "L", "M", "N", "P", "Q")

# Capture current ("before") values for the two row pairs that swap.
$row2 = @{}
$row3 = @{}
$row4 = @{}
$row5 = @{}
foreach ($col in $cols) {
    $row2[$col] = $ws.Range("${col}2").Value2
    $row3[$col] = $ws.Range("${col}3").Value2
    $row4[$col] = $ws.Range("${col}4").Value2
    $row5[$col] = $ws.Range("${col}5").Value2
}

# Row 2 <- old Row 5 ; Row 5 <- old Row 3
# Row 3 <- old Row 4 ; Row 4 <- old Row 2
foreach ($col in $cols) {
    $ws.Range("${col}2").Value = $row5[$col]
    $ws.Range("${col}3").Value = $row4[$col]
    $ws.Range("${col}4").Value = $row2[$col]
    $ws.Range("${col}5").Value = $row3[$col]
}
